$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '88.194.71'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.024.26'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.00'
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '606.62'
$ws.Range("E6").Value = '  -3.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.357'
$ws.Range("E7").Value = '  -7.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.875'
$ws.Range("E8").Value = '  +23.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.021.68'
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.650'
$ws.Range("E11").Value = '  +18.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.184'
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000233'
$ws.Range("E13").Value = '  -5.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.34'
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.279.59'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.606.20'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.42'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.065.30'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.33'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000203'
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.20'
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '418.54'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.98'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.34'
$ws.Range("E25").Value = '  +2.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.91'
$ws.Range("E26").Value = '  +5.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.42'
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.204.51'
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +9.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.160'
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.09'
$ws.Range("E32").Value = '  +0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '496.69'
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.49'
$ws.Range("E34").Value = '  -8.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.52'
$ws.Range("E35").Value = '  -2.88%  '
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.27'
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("E38").Value = '  -2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.16'
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.129'
$ws.Range("E40").Value = '  +4.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +11.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.358'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.88'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.79'
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.31'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0680'
$ws.Range("E48").Value = '  +13.46%  '
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.19'
$ws.Range("E49").Value = '  +2.47%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.98'
$ws.Range("E50").Value = '  +2.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '154.19'
$ws.Range("E51").Value = '  -5.60%  '
